# Fruta / hortaliza, semanal
# Two new weekly price records (Angeleno, fecha 45072) are inserted at the
# top of the "Ciruela" data block (rows 145-146), pushing the existing
# records down by two rows. The table's dimension grows from T230 to T232.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row of the block (145),
# shifting everything below (old rows 145-230) down to 147-232.
$ws.Rows("145:146").Insert()

# --- Row 145: Ciruela / Angeleno / Primera ---
$ws.Range("A145").Value = 5
$ws.Range("B145").Value = "Macroferia Regional de Talca"
$ws.Range("C145").Value = "Maule"
$ws.Range("D145").Value = 45072
$ws.Range("E145").Value = 7
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100103
$ws.Range("H145").Value = "Frutos de hueso (carozo)"
$ws.Range("I145").Value = 100103002
$ws.Range("J145").Value = "Ciruela"
$ws.Range("K145").Value = "Angeleno"
$ws.Range("L145").Value = "Primera"
$ws.Range("M145").Value = 230
$ws.Range("N145").Value = 10000
$ws.Range("O145").Value = 10000
$ws.Range("P145").Value = 10000
$ws.Range("Q145").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R145").Value = "Provincia de Curicó"
$ws.Range("S145").Value = 556
$ws.Range("T145").Value = 18

# --- Row 146: Ciruela / Angeleno / Segunda ---
$ws.Range("A146").Value = 5
$ws.Range("B146").Value = "Macroferia Regional de Talca"
$ws.Range("C146").Value = "Maule"
$ws.Range("D146").Value = 45072
$ws.Range("E146").Value = 7
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100103
$ws.Range("H146").Value = "Frutos de hueso (carozo)"
$ws.Range("I146").Value = 100103002
$ws.Range("J146").Value = "Ciruela"
$ws.Range("K146").Value = "Angeleno"
$ws.Range("L146").Value = "Segunda"
$ws.Range("M146").Value = 200
$ws.Range("N146").Value = 8000
$ws.Range("O146").Value = 8000
$ws.Range("P146").Value = 8000
$ws.Range("Q146").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R146").Value = "Provincia de Curicó"
$ws.Range("S146").Value = 444
$ws.Range("T146").Value = 18
